# Update the "score" sheet's high-score table.
# Three new scores were added to the game (Steven:490, Blaze it:420, Steven:180).
# The table is re-sorted descending by score and truncated to the top 10 rows,
# which drops the two lowest former entries ("Steven N":160 and "NA":80).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("score")

$names  = @("Steven N", "Steven Neveadomi", "Steven", "Blaze it", "Tori", "Steven", "TORI IS MATLAB KING", "Shaleen", "Steven", "Yo Dawg Crilla")
$scores = @(740, 700, 490, 420, 300, 280, 260, 200, 180, 170)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $scores[$i]
}
